$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (B1:F1 unchanged text, G1:I1 new)
$ws.Range("G1").Value = "1 case Spinach Fettuccine"
$ws.Range("H1").Value = "2 egg papp"
$ws.Range("I1").Value = "5 Spinach Linguine"

# Row 2 (gingoso2@gmail.com) - all become "1 P1RAVI", plus new columns
$ws.Range("B2:F2").Value = "1 P1RAVI"
$ws.Range("G2").Value = "1 P2FETT"
$ws.Range("H2").Value = "1 P1PAPP"
$ws.Range("I2").Value = "5 P2LING"

# Row 3 (garrett@gmail.com)
$ws.Range("B3").Value = "1 P1RAVI"
$ws.Range("C3").Value = "1 P1FETT"
$ws.Range("D3").Value = "1 P1LING"
$ws.Range("E3").Value = "1 A2PAPP"
$ws.Range("F3").Value = "1 P1LING"

# Row 4 (DONOVON@gmail.com)
$ws.Range("B4").Value = "1 P1RAVI"
$ws.Range("C4").Value = "2 P1FETT"
$ws.Range("D4").Value = "1 P1LING"
$ws.Range("E4").Value = "1 A3PAPP"
$ws.Range("F4").Value = "1 P1LING"

$ws.Range("J5").Select()
